$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) - temporarily force text format so Excel
# doesn't auto-convert these month/year strings into date serial numbers,
# then restore the default "Normal" style so no formatting change is left
# behind (matches the source edit, which only touched the text content).
$ws.Range("A1").NumberFormat = "@"
$ws.Range("A1").Value = "July 2024"
$ws.Range("A1").Style = "Normal"

$ws.Range("G1").NumberFormat = "@"
$ws.Range("G1").Value = "August 2024"
$ws.Range("G1").Style = "Normal"

# Update data values (row 2)
$ws.Range("A2").Value = 1.549
$ws.Range("B2").Value = -0.392
$ws.Range("C2").Value = 0.013
$ws.Range("D2").Value = -0.063
$ws.Range("E2").Value = -0.051
$ws.Range("F2").Value = 0.039
$ws.Range("G2").Value = 1.094
